$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3606.5386
$ws.Range("I86").Value = 1273.25
$ws.Range("J86").Value = 4643.5557
$ws.Range("K86").Value = 1273.25
$ws.Range("L86").Value = 4643.5557
$ws.Range("M86").Value = -150.25
$ws.Range("N86").Value = -6889.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3606.5386
$ws.Range("I89").Value = 1273.25
$ws.Range("J89").Value = 4643.5557
$ws.Range("K89").Value = 6366.25
$ws.Range("L89").Value = 23217.7785
$ws.Range("M89").Value = -750.25
$ws.Range("N89").Value = -34449.7785

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 596.8333
$ws.Range("J121").Value = 584.8823
$ws.Range("L121").Value = 1754.6469
$ws.Range("N121").Value = -5248.6469

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 38409612
$ws.Range("J129").Value = 2058618.6
$ws.Range("L129").Value = 6175855.800000001
$ws.Range("N129").Value = -6185855.800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3535.7144
$ws.Range("I132").Value = 3653.8462
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 10961.5386
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -8431.5386
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 15626563
$ws.Range("I137").Value = 1663.4348
$ws.Range("J137").Value = 55556860
$ws.Range("K137").Value = 4990.3044
$ws.Range("L137").Value = 166670580
$ws.Range("M137").Value = -2440.3044
$ws.Range("N137").Value = -166675680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2617.372
$ws.Range("I32").Value = 1868.6024
$ws.Range("J32").Value = 23333.334
$ws.Range("K32").Value = 1868.6024
$ws.Range("L32").Value = 23333.334
$ws.Range("M32").Value = -1581.6024
$ws.Range("N32").Value = -23907.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 975.0357
$ws.Range("I74").Value = 1202.5834
$ws.Range("J74").Value = 804.375
$ws.Range("K74").Value = 1202.5834
$ws.Range("L74").Value = 804.375
$ws.Range("M74").Value = -328.5834
$ws.Range("N74").Value = -2552.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 975.0357
$ws.Range("I77").Value = 1202.5834
$ws.Range("J77").Value = 804.375
$ws.Range("K77").Value = 6012.916999999999
$ws.Range("L77").Value = 4021.875
$ws.Range("M77").Value = -1644.916999999999
$ws.Range("N77").Value = -12757.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 949.13336
$ws.Range("I97").Value = 787.46155
$ws.Range("K97").Value = 787.46155
$ws.Range("M97").Value = -291.46155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2851.4285
$ws.Range("I102").Value = 2609.2307
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 2609.2307
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -987.2307000000001
$ws.Range("N102").Value = -9244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1710.5555
$ws.Range("I94").Value = 1132.6666
$ws.Range("J94").Value = 4600
$ws.Range("K94").Value = 1132.6666
$ws.Range("L94").Value = 4600
$ws.Range("M94").Value = -681.6666
$ws.Range("N94").Value = -5502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 960.5
$ws.Range("I107").Value = 967.44446
$ws.Range("J107").Value = 939.6667
$ws.Range("K107").Value = 967.44446
$ws.Range("L107").Value = 939.6667
$ws.Range("M107").Value = 952.55554
$ws.Range("N107").Value = -4779.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2171.4119
$ws.Range("I134").Value = 1637.6364
$ws.Range("J134").Value = 3150
$ws.Range("K134").Value = 4912.9092
$ws.Range("L134").Value = 9450
$ws.Range("M134").Value = -2377.9092
$ws.Range("N134").Value = -14520

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4254.875
$ws.Range("I76").Value = 3013
$ws.Range("K76").Value = 9039
$ws.Range("M76").Value = -8656

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 4254.875
$ws.Range("I79").Value = 3013
$ws.Range("K79").Value = 9039
$ws.Range("M79").Value = -7713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15246.329
$ws.Range("I131").Value = 92293.63
$ws.Range("J131").Value = 1576.6451
$ws.Range("K131").Value = 276880.89
$ws.Range("L131").Value = 4729.9353
$ws.Range("M131").Value = -271840.89
$ws.Range("N131").Value = -14809.9353

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4811.04
$ws.Range("I137").Value = 2019.1666
$ws.Range("J137").Value = 7388.154
$ws.Range("K137").Value = 6057.4998
$ws.Range("L137").Value = 22164.462
$ws.Range("M137").Value = -957.4997999999996
$ws.Range("N137").Value = -32364.462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3879.7144
$ws.Range("I80").Value = 4025.8333
$ws.Range("J80").Value = 3003
$ws.Range("K80").Value = 4025.8333
$ws.Range("L80").Value = 3003
$ws.Range("M80").Value = -3027.8333
$ws.Range("N80").Value = -4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3879.7144
$ws.Range("I83").Value = 4025.8333
$ws.Range("J83").Value = 3003
$ws.Range("K83").Value = 20129.1665
$ws.Range("L83").Value = 15015
$ws.Range("M83").Value = -15137.1665
$ws.Range("N83").Value = -24999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9078.5
$ws.Range("I113").Value = 15837
$ws.Range("J113").Value = 2320
$ws.Range("K113").Value = 15837
$ws.Range("L113").Value = 2320
$ws.Range("M113").Value = -13667
$ws.Range("N113").Value = -6660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4067.8462
$ws.Range("I132").Value = 3285.5
$ws.Range("J132").Value = 5319.6
$ws.Range("K132").Value = 9856.5
$ws.Range("L132").Value = 15958.8
$ws.Range("M132").Value = -7326.5
$ws.Range("N132").Value = -21018.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1735.909
$ws.Range("I82").Value = 2148.6667
$ws.Range("J82").Value = 1240.6
$ws.Range("K82").Value = 2148.6667
$ws.Range("L82").Value = 1240.6
$ws.Range("M82").Value = -1787.6667
$ws.Range("N82").Value = -1962.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1735.909
$ws.Range("I85").Value = 2148.6667
$ws.Range("J85").Value = 1240.6
$ws.Range("K85").Value = 2148.6667
$ws.Range("L85").Value = 1240.6
$ws.Range("M85").Value = -900.6667000000002
$ws.Range("N85").Value = -3736.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4569.5713
$ws.Range("I136").Value = 1724.8572
$ws.Range("J136").Value = 7414.2856
$ws.Range("K136").Value = 5174.571599999999
$ws.Range("L136").Value = 22242.8568
$ws.Range("M136").Value = -2624.571599999999
$ws.Range("N136").Value = -27342.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 49920
$ws.Range("J127").Value = 49920
$ws.Range("L127").Value = 49920
$ws.Range("N127").Value = -59840
